# Edit script: renumber "5.0 General Requirements" -> "4.0 General Requirements",
# collapse the old 5.1/5.2 sub-sections into a single pointer paragraph,
# renumber "4.0 Actions" -> "5.0 Actions", update the two "(see 5.x requirements)"
# cross references to the new "General Requirements" / "GR" wording, and restore
# a lastRenderedPageBreak marker that moved along with the page-break point.

$d = $word.ActiveDocument

$xmlNs = 'xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"'
$pkgOpen = '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document ' + $xmlNs + '><w:body>'
$pkgClose = '</w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'

function Insert-RunXml($range, [string]$text, [bool]$preserve) {
    $space = ""
    if ($preserve) { $space = ' xml:space="preserve"' }
    $escaped = $text -replace '&','&amp;' -replace '<','&lt;' -replace '>','&gt;'
    $frag = $pkgOpen + '<w:p><w:r><w:t' + $space + '>' + $escaped + '</w:t></w:r></w:p>' + $pkgClose
    $range.InsertXML($frag)
}

function Insert-BreakRunXml($range) {
    $frag = $pkgOpen + '<w:p><w:r><w:lastRenderedPageBreak/></w:r></w:p>' + $pkgClose
    $range.InsertXML($frag)
}

# ---------------------------------------------------------------------------
# 1) "5.0 General Requirements" heading -> split into "4" + ".0 General Requirements"
# ---------------------------------------------------------------------------
$p = $d.Paragraphs.Item(19)
if ($p.Range.Text.TrimEnd() -ne "5.0 General Requirements") {
    throw "Unexpected paragraph 19 text: $($p.Range.Text)"
}
$digit = $d.Range($p.Range.Start, $p.Range.Start + 1)
$digit.Delete()
$insPt = $d.Range($p.Range.Start, $p.Range.Start)
Insert-RunXml $insPt "4" $false

# ---------------------------------------------------------------------------
# 2) Remove old "5.1 User Session Management" ... "5.2 Follow-Up List" content
#    (paragraphs 20 through 30) and replace with a single pointer paragraph.
# ---------------------------------------------------------------------------
$pStart = $d.Paragraphs.Item(20)
$pEnd = $d.Paragraphs.Item(30)
if ($pStart.Range.Text.TrimEnd() -ne "5.1 User Session Management") {
    throw "Unexpected paragraph 20 text: $($pStart.Range.Text)"
}
if ($pEnd.Range.Text.TrimEnd() -ne "An Officer should be able to add or move items from a list to another as progress is made or as needed.") {
    throw "Unexpected paragraph 30 text: $($pEnd.Range.Text)"
}
$killRange = $d.Range($pStart.Range.Start, $pEnd.Range.End)
$killRange.Delete()

# Insert the replacement paragraph right after the (now renumbered) heading.
$p19 = $d.Paragraphs.Item(19)
$afterHeading = $d.Range($p19.Range.End - 1, $p19.Range.End - 1)
Insert-RunXml $afterHeading '(see "General Requirements" documentation)' $false

# ---------------------------------------------------------------------------
# 3) "4.0 Actions" heading -> split into "5" + ".0 Actions"
# ---------------------------------------------------------------------------
$pActions = $d.Paragraphs.Item(21)
if ($pActions.Range.Text.TrimEnd() -ne "4.0 Actions") {
    throw "Unexpected paragraph 21 text: $($pActions.Range.Text)"
}
$digit2 = $d.Range($pActions.Range.Start, $pActions.Range.Start + 1)
$digit2.Delete()
$insPt2 = $d.Range($pActions.Range.Start, $pActions.Range.Start)
Insert-RunXml $insPt2 "5" $false

# ---------------------------------------------------------------------------
# 4) "(see 5.1 requirements)" -> "(see " + ""General Requirements" [GR]" + ")"
# ---------------------------------------------------------------------------
$find = $d.Content
$found = $find.Find.Execute("(see 5.1 requirements)", $true, $true, $false, $false, $false, $true, 1, $false, "", 0)
if (-not $found) {
    throw "Could not find '(see 5.1 requirements)'"
}
$target = $d.Range($find.Start, $find.End)
$target.Delete()
$insPt3 = $d.Range($target.Start, $target.Start)
Insert-RunXml $insPt3 "(see " $true
$insPt3b = $d.Range($insPt3.End, $insPt3.End)
Insert-RunXml $insPt3b '"General Requirements" [GR]' $false
$insPt3c = $d.Range($insPt3b.End, $insPt3b.End)
Insert-RunXml $insPt3c ")" $false

# ---------------------------------------------------------------------------
# 5) " (see 5.2 requirements)" -> " (" + "see GR" + ")"
# ---------------------------------------------------------------------------
$find2 = $d.Content
$found2 = $find2.Find.Execute(" (see 5.2 requirements)", $true, $true, $false, $false, $false, $true, 1, $false, "", 0)
if (-not $found2) {
    throw "Could not find ' (see 5.2 requirements)'"
}
$target2 = $d.Range($find2.Start, $find2.End)
$target2.Delete()
$insPt4 = $d.Range($target2.Start, $target2.Start)
Insert-RunXml $insPt4 " (" $true
$insPt4b = $d.Range($insPt4.End, $insPt4.End)
Insert-RunXml $insPt4b "see GR" $false
$insPt4c = $d.Range($insPt4b.End, $insPt4b.End)
Insert-RunXml $insPt4c ")" $false

# ---------------------------------------------------------------------------
# 6) Add <w:lastRenderedPageBreak/> before "Cases near completion"
# ---------------------------------------------------------------------------
$find3 = $d.Content
$found3 = $find3.Find.Execute("Cases near completion", $true, $true, $false, $false, $false, $true, 1, $false, "", 0)
if (-not $found3) {
    throw "Could not find 'Cases near completion'"
}
$breakPt = $d.Range($find3.Start, $find3.Start)
Insert-BreakRunXml $breakPt
